$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.303.88'
$ws.Range("E2").Value = '  +3.70%  '
$ws.Range("D3").Value = '2.254.91'
$ws.Range("E3").Value = '  +2.65%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.62'
$ws.Range("E5").Value = '  +3.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.60'
$ws.Range("E6").Value = '  +5.17%  '
$ws.Range("E7").Value = '  +2.39%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +3.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.91'
$ws.Range("E10").Value = '  +8.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '31.98'
$ws.Range("E11").Value = '  +7.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0792'
$ws.Range("E12").Value = '  +2.26%  '
$ws.Range("E13").Value = '  +3.06%  '
$ws.Range("E14").Value = '  +3.18%  '
$ws.Range("D15").Value = '2.604.65'
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.13'
$ws.Range("E16").Value = '  +3.62%  '
$ws.Range("D17").Value = '2.213.30'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.750'
$ws.Range("E18").Value = '  +4.21%  '
$ws.Range("D19").Value = '41.271.29'
$ws.Range("E19").Value = '  +3.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.17'
$ws.Range("E20").Value = '  +8.67%  '
$ws.Range("D21").Value = '0.0₃0904'
$ws.Range("E21").Value = '  +2.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.89'
$ws.Range("E22").Value = '  +2.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.85'
$ws.Range("E23").Value = '  +2.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '239.73'
$ws.Range("E24").Value = '  +1.54%  '
$ws.Range("E25").Value = '  +4.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("E27").Value = '  +3.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.73'
$ws.Range("E28").Value = '  +5.89%  '
$ws.Range("E29").Value = '  +1.91%  '
$ws.Range("E30").Value = '  +5.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.42'
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.69'
$ws.Range("E32").Value = '  +8.11%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  +6.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0736'
$ws.Range("E35").Value = '  +4.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.05'
$ws.Range("E36").Value = '  +8.57%  '
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.58'
$ws.Range("E38").Value = '  +9.26%  '
$ws.Range("E39").Value = '  +3.01%  '
$ws.Range("E40").Value = '  +6.40%  '
$ws.Range("E41").Value = '  +6.14%  '
$ws.Range("E42").Value = '  +7.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.39'
$ws.Range("E43").Value = '  +17.95%  '
$ws.Range("D44").Value = '2.063.98'
$ws.Range("E44").Value = '  -2.27%  '
$ws.Range("E45").Value = '  +3.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.13'
$ws.Range("E46").Value = '  +5.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("E47").Value = '  +12.60%  '
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("D49").Value = '2.475.86'
$ws.Range("E49").Value = '  +2.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.51'
$ws.Range("E50").Value = '  +2.79%  '
$ws.Range("E51").Value = '  +3.56%  '
